$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-9 from 2023-09-14 (45183)
# to 2023-09-15 (45184), keeping the existing date number format.
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45184
}
